# Row 11 ("Alximik") loses its Russian flag -- the roster's convention is
# to show :skull_crossbones: in the "Flag" column whenever a member has no
# SWGOHGG username on file (see rows 7, 16, 18, 23, 33 for the same
# pattern). Replace the flag and clear the now-stale SWGOHGG handle,
# leaving the cursor on the last cell that was touched (C11).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B11").Value = ":skull_crossbones:"
$ws.Range("C11").ClearContents()

$ws.Range("C11").Select()
